$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update scalar values that changed on this refreshed "estado de cuenta"
$ws.Range("E11").Value = 248202
$ws.Range("C13").Value = 9
$ws.Range("F13").Value = 1
$ws.Range("G16").Value = 689455
$ws.Range("G23").Value = 689455

# 2) Row 24 (last worker row) takes on the closing/bottom-border formatting
#    that the totals row (25) used to have, then the totals row is removed.
$ws.Range("B25:J25").Copy()
$ws.Range("B24:J24").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Rows("25:25").Delete()
